$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1597
$ws.Range("I19").Value = 113.333336
$ws.Range("J19").Value = 3377.4
$ws.Range("K19").Value = 113.333336
$ws.Range("L19").Value = 3377.4
$ws.Range("M19").Value = 61.666664
$ws.Range("N19").Value = -3727.4

$ws.Range("H20").Value = 13094.333
$ws.Range("I20").Value = 1635.5
$ws.Range("J20").Value = 36012
$ws.Range("K20").Value = 1635.5
$ws.Range("L20").Value = 36012
$ws.Range("M20").Value = -1405.5
$ws.Range("N20").Value = -36472

$ws.Range("H35").Value = 13094.333
$ws.Range("I35").Value = 1635.5
$ws.Range("J35").Value = 36012
$ws.Range("K35").Value = 1635.5
$ws.Range("L35").Value = 36012
$ws.Range("M35").Value = -1256.5
$ws.Range("N35").Value = -36770

$ws.Range("H61").Value = 299.3
$ws.Range("I61").Value = 196.625
$ws.Range("K61").Value = 589.875
$ws.Range("M61").Value = -417.875

$ws.Range("H98").Value = 863.4545000000001
$ws.Range("I98").Value = 666.5
$ws.Range("J98").Value = 1749.75
$ws.Range("K98").Value = 666.5
$ws.Range("L98").Value = 1749.75
$ws.Range("M98").Value = 831.5
$ws.Range("N98").Value = -4745.75

$ws.Range("H113").Value = 59725.44
$ws.Range("I113").Value = 97755.734
$ws.Range("J113").Value = 2680
$ws.Range("K113").Value = 97755.734
$ws.Range("L113").Value = 2680
$ws.Range("M113").Value = -94501.734
$ws.Range("N113").Value = -9188

$ws.Range("H116").Value = 3073.182
$ws.Range("I116").Value = 2133.3333
$ws.Range("J116").Value = 3425.625
$ws.Range("K116").Value = 2133.3333
$ws.Range("L116").Value = 3425.625
$ws.Range("M116").Value = 1308.6667
$ws.Range("N116").Value = -10309.625

$ws.Range("H122").Value = 863.4545000000001
$ws.Range("I122").Value = 666.5
$ws.Range("J122").Value = 1749.75
$ws.Range("K122").Value = 1999.5
$ws.Range("L122").Value = 5249.25
$ws.Range("M122").Value = 450.5
$ws.Range("N122").Value = -10149.25

$ws.Range("H137").Value = 6266405
$ws.Range("I137").Value = 891.7059
$ws.Range("K137").Value = 2675.1177
$ws.Range("M137").Value = -125.1177000000002

$ws.Range("H141").Value = 2874.4736
$ws.Range("I141").Value = 1634.5834
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 4903.7502
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = 276.2497999999996
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 12196857
$ws.Range("I74").Value = 14706908
$ws.Range("K74").Value = 14706908
$ws.Range("M74").Value = -14706034

$ws.Range("H77").Value = 12196857
$ws.Range("I77").Value = 14706908
$ws.Range("K77").Value = 73534540
$ws.Range("M77").Value = -73530172

$ws.Range("H102").Value = 1272.8572
$ws.Range("I102").Value = 1282
$ws.Range("K102").Value = 1282
$ws.Range("M102").Value = 340

$ws.Range("H132").Value = 5815946.5
$ws.Range("I132").Value = 7354763.5
$ws.Range("J132").Value = 2638
$ws.Range("K132").Value = 22064290.5
$ws.Range("L132").Value = 7914
$ws.Range("M132").Value = -22061760.5
$ws.Range("N132").Value = -12974

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 19233204
$ws.Range("I86").Value = 2266.5557
$ws.Range("J86").Value = 62502812
$ws.Range("K86").Value = 2266.5557
$ws.Range("L86").Value = 62502812
$ws.Range("M86").Value = -1143.5557
$ws.Range("N86").Value = -62505058

$ws.Range("H89").Value = 19233204
$ws.Range("I89").Value = 2266.5557
$ws.Range("J89").Value = 62502812
$ws.Range("K89").Value = 11332.7785
$ws.Range("L89").Value = 312514060
$ws.Range("M89").Value = -5716.7785
$ws.Range("N89").Value = -312525292

$ws.Range("H135").Value = 46380
$ws.Range("J135").Value = 46380
$ws.Range("L135").Value = 46380
$ws.Range("N135").Value = -56520

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1122.866
$ws.Range("I31").Value = 627.7966300000001
$ws.Range("J31").Value = 1891.5264
$ws.Range("K31").Value = 627.7966300000001
$ws.Range("L31").Value = 1891.5264
$ws.Range("M31").Value = -332.7966300000001
$ws.Range("N31").Value = -2481.5264

$ws.Range("H34").Value = 1122.866
$ws.Range("I34").Value = 627.7966300000001
$ws.Range("J34").Value = 1891.5264
$ws.Range("K34").Value = 627.7966300000001
$ws.Range("L34").Value = 1891.5264
$ws.Range("M34").Value = -425.7966300000001
$ws.Range("N34").Value = -2295.5264

$ws.Range("H132").Value = 12502241
$ws.Range("I132").Value = 15153243
$ws.Range("J132").Value = 4658.143
$ws.Range("K132").Value = 45459729
$ws.Range("L132").Value = 13974.429
$ws.Range("M132").Value = -45457199
$ws.Range("N132").Value = -19034.429

$ws.Range("H140").Value = 39853.332
$ws.Range("J140").Value = 39853.332
$ws.Range("L140").Value = 39853.332
$ws.Range("N140").Value = -50213.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 483.66666
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 483.66666
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1934.66664
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2072.66664

$ws.Range("H113").Value = 684.15625
$ws.Range("I113").Value = 534.1539
$ws.Range("J113").Value = 786.7895
$ws.Range("K113").Value = 1602.4617
$ws.Range("L113").Value = 2360.3685
$ws.Range("M113").Value = 567.5382999999999
$ws.Range("N113").Value = -6700.3685

$ws.Range("H122").Value = 926.3946999999999
$ws.Range("I122").Value = 818.5185
$ws.Range("K122").Value = 7366.6665
$ws.Range("M122").Value = -4916.6665

$ws.Range("H123").Value = 5883.483
$ws.Range("I123").Value = 2158.8
$ws.Range("J123").Value = 7843.8423
$ws.Range("K123").Value = 6476.400000000001
$ws.Range("L123").Value = 23531.5269
$ws.Range("M123").Value = -4026.400000000001
$ws.Range("N123").Value = -28431.5269

$ws.Range("H124").Value = 5851.4443
$ws.Range("I124").Value = 1600
$ws.Range("J124").Value = 7486.615
$ws.Range("K124").Value = 4800
$ws.Range("L124").Value = 22459.845
$ws.Range("M124").Value = 110
$ws.Range("N124").Value = -32279.845

$ws.Range("H125").Value = 5479.231
$ws.Range("I125").Value = 2030
$ws.Range("J125").Value = 5766.6665
$ws.Range("K125").Value = 6090
$ws.Range("L125").Value = 17299.9995
$ws.Range("M125").Value = -1170
$ws.Range("N125").Value = -27139.9995

$ws.Range("H134").Value = 5759.4443
$ws.Range("I134").Value = 3052.8572
$ws.Range("J134").Value = 7481.8184
$ws.Range("K134").Value = 9158.571599999999
$ws.Range("L134").Value = 22445.4552
$ws.Range("M134").Value = -4088.571599999999
$ws.Range("N134").Value = -32585.4552

$ws.Range("H137").Value = 4767.3335
$ws.Range("I137").Value = 3193.077
$ws.Range("J137").Value = 15000
$ws.Range("K137").Value = 9579.231
$ws.Range("L137").Value = 45000
$ws.Range("M137").Value = -4479.231
$ws.Range("N137").Value = -55200

$ws.Range("H139").Value = 2904
$ws.Range("I139").Value = 1583.0769
$ws.Range("J139").Value = 5357.143
$ws.Range("K139").Value = 4749.2307
$ws.Range("L139").Value = 16071.429
$ws.Range("M139").Value = 390.7692999999999
$ws.Range("N139").Value = -26351.429

$ws.Range("H141").Value = 11342
$ws.Range("I141").Value = 6570
$ws.Range("J141").Value = 18500
$ws.Range("K141").Value = 19710
$ws.Range("L141").Value = 55500
$ws.Range("M141").Value = -14530
$ws.Range("N141").Value = -65860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 10475
$ws.Range("J92").Value = 10475
$ws.Range("L92").Value = 10475
$ws.Range("N92").Value = -14219

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 992.8
$ws.Range("I22").Value = 533.75
$ws.Range("J22").Value = 1208.8235
$ws.Range("K22").Value = 533.75
$ws.Range("L22").Value = 1208.8235
$ws.Range("M22").Value = -238.75
$ws.Range("N22").Value = -1798.8235

$ws.Range("H27").Value = 992.8
$ws.Range("I27").Value = 533.75
$ws.Range("J27").Value = 1208.8235
$ws.Range("K27").Value = 533.75
$ws.Range("L27").Value = 1208.8235
$ws.Range("M27").Value = -426.75
$ws.Range("N27").Value = -1422.8235

$ws.Range("H68").Value = 2270
$ws.Range("I68").Value = 2686
$ws.Range("J68").Value = 2010
$ws.Range("K68").Value = 2686
$ws.Range("L68").Value = 2010
$ws.Range("M68").Value = -1937
$ws.Range("N68").Value = -3508

$ws.Range("H71").Value = 2270
$ws.Range("I71").Value = 2686
$ws.Range("J71").Value = 2010
$ws.Range("K71").Value = 13430
$ws.Range("L71").Value = 10050
$ws.Range("M71").Value = -9686
$ws.Range("N71").Value = -17538

$ws.Range("H82").Value = 2353.3333
$ws.Range("I82").Value = 2725
$ws.Range("J82").Value = 2218.182
$ws.Range("K82").Value = 2725
$ws.Range("L82").Value = 2218.182
$ws.Range("M82").Value = -2364
$ws.Range("N82").Value = -2940.182

$ws.Range("H85").Value = 2353.3333
$ws.Range("I85").Value = 2725
$ws.Range("J85").Value = 2218.182
$ws.Range("K85").Value = 2725
$ws.Range("L85").Value = 2218.182
$ws.Range("M85").Value = -1477
$ws.Range("N85").Value = -4714.182

$ws.Range("H122").Value = 6120.3
$ws.Range("I122").Value = 7000.8
$ws.Range("J122").Value = 5239.8
$ws.Range("K122").Value = 21002.4
$ws.Range("L122").Value = 15719.4
$ws.Range("M122").Value = -18552.4
$ws.Range("N122").Value = -20619.4

$ws.Range("H133").Value = 48093
$ws.Range("J133").Value = 48093
$ws.Range("L133").Value = 48093
$ws.Range("N133").Value = -53153

$ws.Range("H136").Value = 17863032
$ws.Range("I136").Value = 20002450
$ws.Range("J136").Value = 34535
$ws.Range("K136").Value = 60007350
$ws.Range("L136").Value = 103605
$ws.Range("M136").Value = -60004800
$ws.Range("N136").Value = -108705

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1818.6666
$ws.Range("I81").Value = 1766.6666
$ws.Range("J81").Value = 1870.6666
$ws.Range("K81").Value = 3533.3332
$ws.Range("L81").Value = 3741.3332
$ws.Range("M81").Value = -2472.3332
$ws.Range("N81").Value = -5863.3332

$ws.Range("H84").Value = 1818.6666
$ws.Range("I84").Value = 1766.6666
$ws.Range("J84").Value = 1870.6666
$ws.Range("K84").Value = 17666.666
$ws.Range("L84").Value = 18706.666
$ws.Range("M84").Value = -12362.666
$ws.Range("N84").Value = -29314.666

$ws.Range("H96").Value = 1484.1562
$ws.Range("I96").Value = 1310.2778
$ws.Range("J96").Value = 1707.7142
$ws.Range("K96").Value = 1310.2778
$ws.Range("L96").Value = 1707.7142
$ws.Range("M96").Value = 62.72219999999993
$ws.Range("N96").Value = -4453.7142

$ws.Range("H136").Value = 920.55
$ws.Range("I136").Value = 970.94116
$ws.Range("K136").Value = 2912.82348
$ws.Range("M136").Value = -362.82348
